# (#33) Alteração nos rótulos da tabela para já transformar a primeira linha
# em cabeçalho automaticamente no Power BI.

$wb = $excel.ActiveWorkbook
$sheets = $wb.Worksheets

$ws1 = $sheets.Item(1)
$ws1.Range("B1").Value = "Ano 2015"
$ws1.Range("C1").Value = "Ano 2030"
$ws1.Range("D1").Value = "Ano 2040"
$ws1.Range("E1").Value = "Ano 2050"

$ws2 = $sheets.Item(2)
$ws2.Range("B1").Value = "Ano 2015"
$ws2.Range("C1").Value = "Ano 2030"
$ws2.Range("D1").Value = "Ano 2040"
$ws2.Range("E1").Value = "Ano 2050"

$ws3 = $sheets.Item(3)
$ws3.Range("B1").Value = "Ano 2015"
$ws3.Range("C1").Value = "Ano 2030"
$ws3.Range("D1").Value = "Ano 2040"
$ws3.Range("E1").Value = "Ano 2050"

$ws4 = $sheets.Item(4)
$ws4.Range("B1").Value = "Intervalo 2015"
$ws4.Range("C1").Value = "Intervalo 2015-2030"
$ws4.Range("D1").Value = "Intervalo 2031-2040"
$ws4.Range("E1").Value = "Intervalo 2041-2050"

$ws5 = $sheets.Item(5)
$ws5.Range("B1").Value = "Ano 2015"
$ws5.Range("C1").Value = "Ano 2030"
$ws5.Range("D1").Value = "Ano 2040"
$ws5.Range("E1").Value = "Ano 2050"

$ws6 = $sheets.Item(6)
$ws6.Range("B1").Value = "Ano 2015"
